$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the "TWB PNPCs" worksheet entirely ---
$pnpc = $wb.Worksheets.Item("TWB PNPCs")
$pnpc.Delete()

# --- Add two new fields to "TWB Episodes": contact entry/exit dates for the
#     primary nominated professional, inserted right after the existing
#     "twb_primary_nominated_professional_consent_date" column (K) and before
#     "twb_previous_suicide_attempts" (formerly L, now shifted to N). ---
$ws = $wb.Worksheets.Item("TWB Episodes")

# Insert two blank columns at L (12) - this pushes the old L/M (previous
# suicide attempts / method of suicide attempt) out to N/O.
$ws.Columns.Item(12).Insert()
$ws.Columns.Item(12).Insert()

# Headers
$ws.Range("L1").Value = "twb_primary_nominated_professional_contact_entry_date"
$ws.Range("M1").Value = "twb_primary_nominated_professional_contact_exit_date"

# Data rows
$ws.Range("L2").Value = 16042020
$ws.Range("M2").Value = 9099999

$ws.Range("L3").Value = 9099999
$ws.Range("M3").Value = 9099999
